$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Cow Milk" (index 2): add "cid" column (M) with per-row company id
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("M1").Value = "cid"
$ws2.Range("M2").Value = 29
$ws2.Range("M3").Value = 30
$ws2.Range("M4").Value = 43
$ws2.Range("M5").Value = 25
$ws2.Range("M6").Value = 25

# ---------------------------------------------------------------------
# Sheet "Flavoured Milk" (index 3): add "cid" column (N)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("N1").Value = "cid"
$ws3.Range("N2").Value = 32
$ws3.Range("N3").Value = 12
$ws3.Range("N4").Value = 42
$ws3.Range("N5").Value = 54

# ---------------------------------------------------------------------
# Sheet "Cone" (index 4): add "cid" column (M) + explicit column width
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("M1").Value = "cid"
$ws4.Range("M2").Value = 210
$ws4.Range("M3").Value = 23
$ws4.Range("M4").Value = 42
$ws4.Range("M5").Value = 35
$ws4.Columns.Item(13).ColumnWidth = 9.33

# ---------------------------------------------------------------------
# Row heights: drop the per-row 15.75pt override (sheet1 "Catalog Structure")
# so the rows fall back to the sheet default.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:F11").EntireRow.AutoFit() | Out-Null
$ws2.Range("A1:M6").EntireRow.AutoFit() | Out-Null
$ws3.Range("A1:N5").EntireRow.AutoFit() | Out-Null
$ws4.Range("A1:M5").EntireRow.AutoFit() | Out-Null

# ---------------------------------------------------------------------
# Selections / active sheet: "Cow Milk" becomes the active tab; "Cone"
# loses tabSelected. Set selections to match the target workbook state.
# ---------------------------------------------------------------------
$ws3.Range("A4").Select()
$ws4.Range("M2").Select()
$ws2.Activate()
$ws2.Range("M5").Select()

Write-Output "done"
